$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "GiaoVien" (Teachers) worksheet at the end of the
#    workbook and populate it with the teacher list.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$gv = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$gv.Name = "GiaoVien"

$gvCodes   = @("GV1", "GV2", "GV3", "GV4", "GV5", "GV6")
$gvNames   = @("Dũng", "Đăng", "Quý", "Hùng", "Minh", "Óc")
$gvEmails  = @("dung@gmail.com", "dang@gmail.com", "quy@gmail.com", "hung@gmail.com", "minh@gmail.com", "oc@gmail.com")

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 1
    $gv.Cells.Item($r, 1).Value = $gvCodes[$i]
    $gv.Cells.Item($r, 2).Value = $gvNames[$i]
    $gv.Cells.Item($r, 3).Value = $gvEmails[$i]
}

# Rows 1-2 pick up the auto "looks like a hyperlink" styling (no
# underline); rows 3-6 use the regular Hyperlink cell style.
$gv.Range("C1:C2").Style = "Hyperlink"
$gv.Range("C1:C2").Font.Underline = $false
$gv.Range("C3:C6").Style = "Hyperlink"

$gv.Columns.Item(3).ColumnWidth = 20.43

$gv.Range("A1:A6").Select()

# ------------------------------------------------------------------
# 2. "Lop" (Class) sheet: add a "Giáo Viên" column (E) assigning a
#    teacher to each class.
# ------------------------------------------------------------------
$lop = $wb.Worksheets.Item("Lop")
$lop.Activate()

$lop.Range("E1").Value = "GV1"
$lop.Range("E2").Value = "GV2"
$lop.Range("E3").Value = "GV3"
$lop.Range("E4").Value = "GV4"
$lop.Range("E5").Value = "GV5"

$lop.Range("E1:E5").Select()

# ------------------------------------------------------------------
# 3. "SinhVien" (Student) sheet: row 6 gets a student code in column
#    E, and its enrollment-date value is normalised to a whole day.
# ------------------------------------------------------------------
$sv = $wb.Worksheets.Item("SinhVien")
$sv.Range("E6").Value = "1234a"
$sv.Range("F6").Value = 44120

# ------------------------------------------------------------------
# 4. Leave the workbook focused back on "Lop" with the new column
#    selected, matching the final author view.
# ------------------------------------------------------------------
$lop.Activate()
$lop.Range("E1:E5").Select()
